# Logbook.docx edit: split the "18.05" closing paragraph into the
# original "Haven't work this Monday." remark, a new "19.05" day
# heading, and the full set of 19.05 log entries (ngSwitch work,
# notes feature, deleteElem bug, Mr.Erhler validation, ...).

$d = $word.ActiveDocument

# The paragraph we are about to split is the very last paragraph in
# the document body; it currently holds the lone sentence about the
# large ngSwitch and also carries the trailing _GoBack bookmark.
$lastIndex = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($lastIndex)
$targetRange = $target.Range

# Namespace boilerplate used to wrap every WordprocessingML fragment we
# feed to Range.InsertXML so it is parsed as real document content
# (runs, proofErr marks, paragraph styles, ...) instead of plain text.
$pkgHead = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function ConvertTo-PackageXml([string]$bodyXml) {
    return $pkgHead + $bodyXml + $pkgTail
}

# The 8 paragraphs that replace the original single paragraph, in
# final document order. The last one reuses the original paragraph
# (and therefore keeps the _GoBack bookmark attached at the very end
# of the story, exactly where it already lives).
$p1 = '<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>Haven’t work this</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>M</w:t></w:r><w:r><w:t>onday</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'

$p2 = '<w:p><w:pPr><w:pStyle w:val="Titre1"/></w:pPr><w:r><w:t>19</w:t></w:r><w:r><w:t>.05</w:t></w:r></w:p>'

$p3 = '<w:p><w:r><w:t xml:space="preserve">I’ve linked all interventions icons to the UI using a large </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ngSwitch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>'

$p4 = '<w:p><w:r><w:t xml:space="preserve">I’ve replaced the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ngSwitch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> directive by a simple scope variable that I forged with a switch before loading it on the page.</w:t></w:r></w:p>'

$p5 = '<w:p><w:r><w:t xml:space="preserve">I’ve added the possibility to take </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>notes,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> it will then delete the element from the list and save it into an array which will be later on submitted to the server.</w:t></w:r></w:p>'

$p6 = '<w:p><w:r><w:t>The same functionalities are to be found on the child state/view details.</w:t></w:r></w:p>'

$p7 = '<w:p><w:r><w:t xml:space="preserve">Problem remains with the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>deleteElem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> function when in details view deleting an element doesn’t work anymore.</w:t></w:r></w:p>'

$p8 = '<w:p><w:r><w:t xml:space="preserve">The main view is almost over, once this problem has been solved and the view validated by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mr.Erhler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> on Thursday, I will then move on to the vitals measurement view.</w:t></w:r></w:p>'

$leadingParagraphs = @($p1, $p2, $p3, $p4, $p5, $p6, $p7)

# Create 7 empty paragraphs immediately before the target paragraph.
# InsertParagraphBefore keeps the (still untouched) target paragraph -
# and its bookmark - as the very last paragraph of the story.
foreach ($unused in $leadingParagraphs) {
    [void]$targetRange.InsertParagraphBefore()
}

# Fill the 7 freshly created (empty) paragraphs in order: they are
# simple Range.InsertXML targets since none of them is the story's
# final paragraph.
$baseIndex = $d.Paragraphs.Count - $leadingParagraphs.Count
for ($i = 0; $i -lt $leadingParagraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($baseIndex + $i)
    $xml = ConvertTo-PackageXml $leadingParagraphs[$i]
    [void]$para.Range.InsertXML($xml)
}

# The 8th (last) paragraph is the original one, which still carries the
# _GoBack bookmark. Replacing its Range outright via InsertXML would
# push the new content into a fresh paragraph *before* it (Word always
# keeps a trailing, now-empty paragraph mark at the very end of a
# story) and strand the bookmark in a stray empty paragraph. Instead,
# delete just the old run text (excluding the final paragraph mark, so
# the bookmark stays put) and insert the new XML at the now-collapsed,
# still-last-paragraph insertion point.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
if ($finalRange.End -gt $finalRange.Start) {
    $oldTextRange = $d.Range($finalRange.Start, $finalRange.End - 1)
    if ($oldTextRange.End -gt $oldTextRange.Start) {
        $oldTextRange.Delete()
    }
}
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($finalPara.Range.Start, $finalPara.Range.Start)
[void]$insertionPoint.InsertXML((ConvertTo-PackageXml $p8))

Write-Host "Split last paragraph into 8 paragraphs for 19.05."
